$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.774.91"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.648.71"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.880.90"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "1.658.42"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "26.811.54"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").Value = "1.285.09"
$ws.Range("E34").Value = "  +3.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("E38").Value = "  +6.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").Value = "1.791.25"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.99%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.33%  "
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
